$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: header "Week 1" -> "Week 3"
$ws.Range("B1").Value = "Week 3"

# Update the active selection to B2, matching the saved workbook state
$ws.Range("B2").Select()
